$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.039690384347496
$ws.Range("D2").Value = 1.043064676187131
$ws.Range("E2").Value = 1.047083747646935
$ws.Range("F2").Value = 1.055334453513963
$ws.Range("I2").Value = 1.041427754092699
$ws.Range("J2").Value = 1.044781034722958
$ws.Range("K2").Value = 1.045839595071605
$ws.Range("L2").Value = 1.049847376264761
$ws.Range("M2").Value = 1.05807519649695
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.040609577373308
$ws.Range("D3").Value = 1.043770930488109
$ws.Range("E3").Value = 1.047935789649966
$ws.Range("F3").Value = 1.056349101702798
$ws.Range("I3").Value = 1.041679943851033
$ws.Range("J3").Value = 1.045345689014042
$ws.Range("K3").Value = 1.04635708247969
$ws.Range("L3").Value = 1.050511086525434
$ws.Range("M3").Value = 1.058902752084843
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.041204687713334
$ws.Range("D4").Value = 1.044228070383623
$ws.Range("E4").Value = 1.048487844858186
$ws.Range("F4").Value = 1.057006714139266
$ws.Range("I4").Value = 1.041841896329369
$ws.Range("J4").Value = 1.04571075354616
$ws.Range("K4").Value = 1.046691403682436
$ws.Range("L4").Value = 1.050940633937387
$ws.Range("M4").Value = 1.059438690554835
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.041454950318723
$ws.Range("D5").Value = 1.044420285625037
$ws.Range("E5").Value = 1.04872010162592
$ws.Range("F5").Value = 1.057283428643119
$ws.Range("I5").Value = 1.041909685801265
$ws.Range("J5").Value = 1.045864152826446
$ws.Range("K5").Value = 1.04683182498823
$ws.Range("L5").Value = 1.051121234383256
$ws.Range("M5").Value = 1.059664106800136
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.04149697506855
$ws.Range("D6").Value = 1.04445256134245
$ws.Range("E6").Value = 1.048759108652024
$ws.Range("F6").Value = 1.05732990509794
$ws.Range("I6").Value = 1.041921050614766
$ws.Range("J6").Value = 1.045889904882781
$ws.Range("K6").Value = 1.046855394848232
$ws.Range("L6").Value = 1.051151559060126
$ws.Range("M6").Value = 1.059701961462942
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.041208031428442
$ws.Range("D7").Value = 1.044230638642651
$ws.Range("E7").Value = 1.048490947607293
$ws.Range("F7").Value = 1.057010410614608
$ws.Range("I7").Value = 1.041842803296893
$ws.Range("J7").Value = 1.045712803566204
$ws.Range("K7").Value = 1.046693280500341
$ws.Range("L7").Value = 1.050943047057147
$ws.Range("M7").Value = 1.059441702155922
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.04000096078652
$ws.Range("D8").Value = 1.043303327061337
$ws.Range("E8").Value = 1.047371547906698
$ws.Range("F8").Value = 1.055677136890551
$ws.Range("I8").Value = 1.041513237355955
$ws.Range("J8").Value = 1.04497192477921
$ws.Range("K8").Value = 1.046014591048243
$ws.Range("L8").Value = 1.050071662551821
$ws.Range("M8").Value = 1.058354777957166
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.037876528091758
$ws.Range("D9").Value = 1.041670464601232
$ws.Range("E9").Value = 1.045404648419297
$ws.Range("F9").Value = 1.053335967298721
$ws.Range("I9").Value = 1.040923097686925
$ws.Range("J9").Value = 1.043664110595815
$ws.Range("K9").Value = 1.044814654877621
$ws.Range("L9").Value = 1.048536849176096
$ws.Range("M9").Value = 1.056443010994521
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.036462033688124
$ws.Range("D10").Value = 1.040582762277525
$ws.Range("E10").Value = 1.044097241174736
$ws.Range("F10").Value = 1.051780794084997
$ws.Range("I10").Value = 1.040523385298083
$ws.Range("J10").Value = 1.042790749141807
$ws.Range("K10").Value = 1.044012066693779
$ws.Range("L10").Value = 1.047514156549238
$ws.Range("M10").Value = 1.055170941308986
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035849979625534
$ws.Range("D11").Value = 1.040111998788371
$ws.Range("E11").Value = 1.043532050627646
$ws.Range("F11").Value = 1.051108731918158
$ws.Range("I11").Value = 1.040348821826855
$ws.Range("J11").Value = 1.042412233057824
$ws.Range("K11").Value = 1.043663925119563
$ws.Range("L11").Value = 1.047071454383241
$ws.Range("M11").Value = 1.054620713970781
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.035622701221376
$ws.Range("D12").Value = 1.039937170454716
$ws.Range("E12").Value = 1.043322253981537
$ws.Range("F12").Value = 1.050859300329804
$ws.Range("I12").Value = 1.040283758576276
$ws.Range("J12").Value = 1.04227158440704
$ws.Range("K12").Value = 1.043534518279036
$ws.Range("L12").Value = 1.046907035672839
$ws.Range("M12").Value = 1.05441642433547
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.035671450221495
$ws.Range("D13").Value = 1.039974670157615
$ws.Range("E13").Value = 1.043367249705812
$ws.Range("F13").Value = 1.050912795082025
$ws.Range("I13").Value = 1.040297724928802
$ws.Range("J13").Value = 1.042301756310341
$ws.Range("K13").Value = 1.043562280613951
$ws.Range("L13").Value = 1.046942303097063
$ws.Range("M13").Value = 1.054460241097606
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035831191375605
$ws.Range("D14").Value = 1.040097546718766
$ws.Range("E14").Value = 1.043514705902337
$ws.Range("F14").Value = 1.051088109693053
$ws.Range("I14").Value = 1.040343448216122
$ws.Range("J14").Value = 1.042400608030209
$ws.Range("K14").Value = 1.043653230172165
$ws.Range("L14").Value = 1.047057863055736
$ws.Range("M14").Value = 1.054603825497956
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035929621909787
$ws.Range("D15").Value = 1.040173259592961
$ws.Range("E15").Value = 1.043605577158775
$ws.Range("F15").Value = 1.051196153643307
$ws.Range("I15").Value = 1.040371590361818
$ws.Range("J15").Value = 1.042461507115794
$ws.Range("K15").Value = 1.043709255092002
$ws.Range("L15").Value = 1.047129066117353
$ws.Range("M15").Value = 1.054692304450727
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036502662867062
$ws.Range("D16").Value = 1.040614010021956
$ws.Range("E16").Value = 1.044134770614501
$ws.Range("F16").Value = 1.051825424931309
$ws.Range("I16").Value = 1.040534939253072
$ws.Range("J16").Value = 1.042815862818024
$ws.Range("K16").Value = 1.044035158812755
$ws.Range("L16").Value = 1.047543540066951
$ws.Range("M16").Value = 1.055207470549033
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036862232453386
$ws.Range("D17").Value = 1.040890540628318
$ws.Range("E17").Value = 1.04446696833425
$ws.Range("F17").Value = 1.052220509271004
$ws.Range("I17").Value = 1.040637006412749
$ws.Range("J17").Value = 1.043038049093456
$ws.Range("K17").Value = 1.044239425448929
$ws.Range("L17").Value = 1.047803564253176
$ws.Range("M17").Value = 1.055530778606474
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.037072004866549
$ws.Range("D18").Value = 1.041051857274865
$ws.Range("E18").Value = 1.044660822849026
$ws.Range("F18").Value = 1.052451084309085
$ws.Range("I18").Value = 1.040696397106258
$ws.Range("J18").Value = 1.043167613260471
$ws.Range("K18").Value = 1.044358511263767
$ws.Range("L18").Value = 1.04795524448632
$ws.Range("M18").Value = 1.055719415313594
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.037143538859233
$ws.Range("D19").Value = 1.041106865608258
$ws.Range("E19").Value = 1.044726937321499
$ws.Range("F19").Value = 1.052529726287539
$ws.Range("I19").Value = 1.040716623465784
$ws.Range("J19").Value = 1.043211785613355
$ws.Range("K19").Value = 1.044399106335839
$ws.Range("L19").Value = 1.048006965620623
$ws.Range("M19").Value = 1.055783745121486
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.036823649734588
$ws.Range("D20").Value = 1.040860869327348
$ws.Range("E20").Value = 1.044431317426811
$ws.Range("F20").Value = 1.052178107105898
$ws.Range("I20").Value = 1.040626070396212
$ws.Range("J20").Value = 1.043014214049778
$ws.Range("K20").Value = 1.044217515715342
$ws.Range("L20").Value = 1.047775664816255
$ws.Range("M20").Value = 1.055496084850819
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035784149765377
$ws.Range("D21").Value = 1.04006136167984
$ws.Range("E21").Value = 1.043471279844989
$ws.Range("F21").Value = 1.051036478317765
$ws.Range("I21").Value = 1.040329989988151
$ws.Range("J21").Value = 1.042371500061194
$ws.Range("K21").Value = 1.043626450314154
$ws.Range("L21").Value = 1.047023832952025
$ws.Range("M21").Value = 1.054561540997145
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.03513095560871
$ws.Range("D22").Value = 1.039558877934955
$ws.Range("E22").Value = 1.042868478399376
$ws.Range("F22").Value = 1.050319861584365
$ws.Range("I22").Value = 1.040142544836716
$ws.Range("J22").Value = 1.041967106074668
$ws.Range("K22").Value = 1.043254294765117
$ws.Range("L22").Value = 1.046551246005844
$ws.Range("M22").Value = 1.053974473361649
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.035477189792415
$ws.Range("D23").Value = 1.039825234826178
$ws.Range("E23").Value = 1.043187957436619
$ws.Range("F23").Value = 1.05069964231641
$ws.Range("I23").Value = 1.040242034926423
$ws.Range("J23").Value = 1.042181510570191
$ws.Range("K23").Value = 1.043451631322228
$ws.Range("L23").Value = 1.046801761624401
$ws.Range("M23").Value = 1.054285639641836
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.036841083476087
$ws.Range("D24").Value = 1.040874276446247
$ws.Range("E24").Value = 1.044447426261797
$ws.Range("F24").Value = 1.052197266420554
$ws.Range("I24").Value = 1.040631012354329
$ws.Range("J24").Value = 1.043024984183099
$ws.Range("K24").Value = 1.044227415964651
$ws.Range("L24").Value = 1.047788271332266
$ws.Range("M24").Value = 1.055511761290866
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.038425433012668
$ws.Range("D25").Value = 1.042092450723106
$ws.Range("E25").Value = 1.045912464504228
$ws.Range("F25").Value = 1.053940232979099
$ws.Range("I25").Value = 1.041076773494662
$ws.Range("J25").Value = 1.044002477745125
$ws.Range("K25").Value = 1.045125335051925
$ws.Range("L25").Value = 1.048933548678691
$ws.Range("M25").Value = 1.056936823114824